# Fix a key bug in the dialogue system log: add the missing log entry
# for row 6 (date 45618 already present in column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Impoved dialogue system and added decision support"
$ws.Range("C6").Value = 6
